$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UMT-L2Manager Scenarios cred")

# New values for the "L2 manager" test cases (entered in authoring order so the
# shared-string table gets new entries in this sequence)
$ws.Range("A6").Value = "Roles"
$ws.Range("A9").Value = "TrainingWorkOrderApprover"
$ws.Range("A10").Value = "RGTST01"
$ws.Range("A7").Value = "Upstream|AFT|Alaska|Reporting - Standard User"

# A6/A9 reuse the bold+green header style already used by A1 ("UserName")
[void]$ws.Range("A1").Copy()
[void]$ws.Range("A6").PasteSpecial(-4122)
[void]$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# A7 gets its own look: Times New Roman, vertically centered
$ws.Range("A7").Font.Name = "Times New Roman"
$ws.Range("A7").Font.Family = 1
$ws.Range("A7").VerticalAlignment = -4108

# Widen column A (stored width "37") and move the visible selection to the
# newly added row
$ws.Columns("A").ColumnWidth = 36.17
[void]$ws.Range("A7").Select()

# Sheet1 view no longer pins a frozen top row
$ws1 = $wb.Worksheets.Item("Roles")
$ws1.Application.ActiveWindow.ScrollRow = 1
